$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "wallet" example row (row 26), mirroring the pattern used by the
# existing wallet-related rows 24-25: column B holds the link (styled like a
# hyperlink), column C repeats the "Inspiración para Wallet" description.
$newUrl = "https://themeforest.net/item/crypto-admin-responsive-bootstrap-4-admin-html-templates/21604673"

$ws.Range("B26").Value = $newUrl
$ws.Range("C26").Value = $ws.Range("C25").Value2

# Turn B26 into a real hyperlink, same as the other link cells above it.
$ws.Hyperlinks.Add($ws.Range("B26"), $newUrl) | Out-Null

# Hyperlinks.Add() stamps a fresh direct-format style on the cell; put it
# back in line with the rest of the column (the "Hipervínculo" cell style).
$ws.Range("B26").Style = $ws.Range("B25").Style

# Match the selection left behind by the edit (B2:C28 -> C25:C26 selected).
$ws.Range("C25:C26").Select() | Out-Null
